$d = $word.ActiveDocument

# --- First paragraph text edits -------------------------------------------------
# Replace the old topic ID placeholder text with the new one, and at the same time
# eat the following run's single trailing space (the old paragraph had two runs:
# the ID text, then a run containing just " "; the new paragraph has a single run
# with just the updated ID text and no trailing space run).
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_13__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFICC_PGI_5301__ID**", 2)

# --- First paragraph formatting edits -------------------------------------------
$p1 = $d.Paragraphs(1)

# w:ind w:left="120" -> w:ind w:left="225"  (225 twips = 11.25 pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (top/left/bottom/right) with 5-twip spacing and no line,
# matching <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
